# Update the "Metadata" worksheet (sheet 1) with the new ValueSet metadata:
#  - Version bump 0.1.6 -> 0.1.7
#  - Status active -> draft
#  - Date refreshed
#  - Contact details replaced (org + named contact), new Jurisdiction row added
#  - Existing Description/Purpose/Copyright/Immutable rows shifted down one row
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Extend formatting down to the new last row (16) by copying the format
# already used by row 15 so the new row matches the sheet's existing style.
$ws1.Range("A15:B15").Copy() | Out-Null
$ws1.Range("A16:B16").PasteSpecial(-4122) | Out-Null

# Version
$ws1.Range("B3").Value = "0.1.7"

# Status
$ws1.Range("B6").Value = "draft"

# Date
$ws1.Range("B8").Value = "2024-08-23T10:17:11-05:00"

# Contact details: row 10 keeps "Contact" label but gets the org contact
# string (with URL); row 11 keeps "Contact" label and becomes the named
# contact person.
$ws1.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"
$ws1.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# New Jurisdiction row (was previously the Description row; Description and
# everything below it shifts down by one row).
$ws1.Range("A12").Value = "Jurisdiction"
$ws1.Range("B12").Value = ""

$ws1.Range("A13").Value = "Description"
$ws1.Range("B13").Value = "RxNorm codes for Cytarabine"

$ws1.Range("A14").Value = "Purpose"
$ws1.Range("B14").Value = ""

$ws1.Range("A15").Value = "Copyright"
$ws1.Range("B15").Value = ""

$ws1.Range("A16").Value = "Immutable"
$ws1.Range("B16").Value = "BooleanType[null]"
